# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" fund-holdings sheet (same shape as the existing
#   quarterly sheets) right before the "总计" (total) summary sheet.
# - Rebuilds "总计" with a new first data row for 2022-Q1 and the previous
#   rows shifted down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Work out sheet identity up front (by name, not by captured object
#    reference -- indices/handles shift as sheets are added/moved).
# ---------------------------------------------------------------------
$originalActiveSheetName = $wb.ActiveSheet.Name
$totalSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name
$lastQuarterName = $wb.Worksheets.Item($wb.Worksheets.Count - 1).Name

# Snapshot the old "总计" rows (2021-Q4 .. 2020-Q4) before we remove the sheet.
# NB: `.Value` as a getter is unreliable on this host -- use `.Value2` to read.
$oldTotal = $wb.Worksheets.Item($totalSheetName)
$oldRows = @()
for ($r = 2; $r -le 6; $r++) {
    $oldRows += , @($oldTotal.Cells.Item($r, 2).Value2, $oldTotal.Cells.Item($r, 3).Value2, $oldTotal.Cells.Item($r, 4).Value2)
}

# ---------------------------------------------------------------------
# 2) Delete the old "总计" sheet so its sheetId is freed, then recreate
#    the "2022-Q1" sheet (reuses the freed id) followed by a fresh "总计"
#    sheet (gets the next id) -- this reproduces the sheetId renumbering
#    seen in the target workbook (2022-Q1 -> 6, 总计 -> 7).
# ---------------------------------------------------------------------
$oldTotal.Delete()

$newQuarter = $wb.Worksheets.Add()
$newQuarter.Name = "2022-Q1"

$newTotal = $wb.Worksheets.Add()
$newTotal.Name = $totalSheetName

# ---------------------------------------------------------------------
# 3) Fix up tab order: 2020-Q4, 2021-Q1 .. 2021-Q4, 2022-Q1, 总计
# ---------------------------------------------------------------------
$lastQuarter = $wb.Worksheets.Item($lastQuarterName)
$quarterSheet = $wb.Worksheets.Item("2022-Q1")
$lastQuarter.Move($quarterSheet)

$quarterSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$totalSheet.Move($null, $quarterSheet)

# ---------------------------------------------------------------------
# 4) Populate "2022-Q1" -- copy header/column-A formatting from the
#    "2021-Q4" sheet (bold, bordered, centered) then write the new values.
# ---------------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item("2022-Q1")
$template = $wb.Worksheets.Item($lastQuarterName)
$template.Range("A1:H5").Copy()
$quarterSheet.Range("A1").PasteSpecial(-4122)

$quarterSheet.Range("B1").Value = "基金代码"
$quarterSheet.Range("C1").Value = "基金名称"
$quarterSheet.Range("D1").Value = "基金规模"
$quarterSheet.Range("E1").Value = "股票总仓位"
$quarterSheet.Range("F1").Value = "仓位占比"
$quarterSheet.Range("G1").Value = "持有市值(亿元)"
$quarterSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @("167301", "方正富邦中证保险主题指数（LOF）", "54.05", "93.05", "3.35", "1.8107", 7),
    @("257040", "国联安红利混合", "1.10", "72.31", "4.98", "0.0548", 5),
    @("005248", "新华沪深300指数增强A", "2.26", "94.10", "1.80", "0.0407", 10),
    @("510760", "国泰上证综合ETF", "2.22", "95.43", "1.35", "0.0300", 10)
)

$quarterSheet.Range("B2:G5").NumberFormat = "@"

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $row = $i + 2
    $data = $fundRows[$i]
    $quarterSheet.Cells.Item($row, 1).Value = $i
    $quarterSheet.Cells.Item($row, 2).Value = $data[0]
    $quarterSheet.Cells.Item($row, 3).Value = $data[1]
    $quarterSheet.Cells.Item($row, 4).Value = $data[2]
    $quarterSheet.Cells.Item($row, 5).Value = $data[3]
    $quarterSheet.Cells.Item($row, 6).Value = $data[4]
    $quarterSheet.Cells.Item($row, 7).Value = $data[5]
    $quarterSheet.Cells.Item($row, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 5) Populate "总计" -- same header as before, plus the new 2022-Q1 row
#    inserted at the top and the previous rows shifted down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$template.Range("A1:D1").Copy()
$totalSheet.Range("A1").PasteSpecial(-4122)
$template.Range("A2:A2").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(, @("2022-Q1", 4, 1.94))
$totalRows += $oldRows

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $row = $i + 2
    $data = $totalRows[$i]
    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $data[0]
    $totalSheet.Cells.Item($row, 3).Value = $data[1]
    $totalSheet.Cells.Item($row, 4).Value = $data[2]
}
